$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '54.651.93'
$ws.Range("E2").Value = '  -1.33%  '
$ws.Range("D3").Value = '2.295.33'
$ws.Range("E3").Value = '  -1.76%  '
$ws.Range("D4").Value = '0.995'
$ws.Range("E4").Value = '  -0.73%  '
$ws.Range("D5").Value = '495.19'
$ws.Range("E5").Value = '  -0.71%  '
$ws.Range("D6").Value = '127.28'
$ws.Range("E6").Value = '  -1.70%  '
$ws.Range("D7").Value = '0.995'
$ws.Range("E7").Value = '  -0.51%  '
$ws.Range("D8").Value = '0.528'
$ws.Range("E8").Value = '  -0.67%  '
$ws.Range("D9").Value = '2.293.70'
$ws.Range("E9").Value = '  -1.97%  '
$ws.Range("E10").Value = '  +0.38%  '
$ws.Range("D11").Value = '0.152'
$ws.Range("E11").Value = '  +2.03%  '
$ws.Range("E12").Value = '  +2.49%  '
$ws.Range("E13").Value = '  -2.25%  '
$ws.Range("D14").Value = '2.695.96'
$ws.Range("E14").Value = '  -2.24%  '
$ws.Range("D15").Value = '21.73'
$ws.Range("E15").Value = '  +1.20%  '
$ws.Range("D16").Value = '54.259.38'
$ws.Range("E16").Value = '  -2.03%  '
$ws.Range("E17").Value = '  +0.20%  '
$ws.Range("D18").Value = '2.299.04'
$ws.Range("E18").Value = '  -2.18%  '
$ws.Range("E19").Value = '  +3.05%  '
$ws.Range("E20").Value = '  +2.25%  '
$ws.Range("D21").Value = '305.22'
$ws.Range("E21").Value = '  -0.47%  '
$ws.Range("D22").Value = '6.44'
$ws.Range("E22").Value = '  +4.50%  '
$ws.Range("D23").Value = '0.999'
$ws.Range("E23").Value = '  -0.01%  '
$ws.Range("E24").Value = '  -2.53%  '
$ws.Range("D25").Value = '62.91'
$ws.Range("E25").Value = '  -2.78%  '
$ws.Range("E26").Value = '  +0.48%  '
$ws.Range("B27").Value = 'Kaspa'
$ws.Range("C27").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D27").Value = '0.152'
$ws.Range("E27").Value = '  +5.84%  '
$ws.Range("B28").Value = 'Polygon'
$ws.Range("C28").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D28").Value = '0.374'
$ws.Range("E28").Value = '  +0.61%  '
$ws.Range("D29").Value = '2.382.21'
$ws.Range("E29").Value = '  -3.43%  '
$ws.Range("D30").Value = '7.12'
$ws.Range("E30").Value = '  +0.19%  '
$ws.Range("D31").Value = '171.87'
$ws.Range("E31").Value = '  +2.15%  '
$ws.Range("E32").Value = '  -1.65%  '
$ws.Range("D33").Value = '0.0₃0685'
$ws.Range("E33").Value = '  -1.98%  '
$ws.Range("D34").Value = '5.89'
$ws.Range("E34").Value = '  +3.12%  '
$ws.Range("E35").Value = '  -0.27%  '
$ws.Range("D36").Value = '0.999'
$ws.Range("E36").Value = '  -0.16%  '
$ws.Range("E37").Value = '  +1.11%  '
$ws.Range("D38").Value = '17.61'
$ws.Range("E38").Value = '  +0.37%  '
$ws.Range("E39").Value = '  +2.79%  '
$ws.Range("D40").Value = '0.867'
$ws.Range("E40").Value = '  +2.70%  '
$ws.Range("D41").Value = '3.66'
$ws.Range("E41").Value = '  +0.64%  '
$ws.Range("D42").Value = '35.55'
$ws.Range("E42").Value = '  -1.53%  '
$ws.Range("E43").Value = '  +0.43%  '
$ws.Range("E44").Value = '  +1.54%  '
$ws.Range("E45").Value = '  -0.13%  '
$ws.Range("D46").Value = '128.70'
$ws.Range("E46").Value = '  +4.35%  '
$ws.Range("D47").Value = '4.88'
$ws.Range("E47").Value = '  +3.21%  '
$ws.Range("D48").Value = '0.0895'
$ws.Range("E48").Value = '  +0.87%  '
$ws.Range("D49").Value = '0.552'
$ws.Range("E49").Value = '  +0.54%  '
$ws.Range("D50").Value = '243.29'
$ws.Range("E50").Value = '  +1.90%  '
$ws.Range("D51").Value = '0.0483'
$ws.Range("E51").Value = '  +1.65%  '
